$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "relation"
$ws.Range("B1").Value = "count"

$wb.Save()
